# Added functional healthbar slider
# - Merge the "Player dash " task note into the "Add, test, finalize dash animation"
#   note so that row 7 / column A reads "Add, test, finalize dash and animation"
# - Clear the now-redundant note that used to live in A8
# - Leave the cursor/selection on B18, matching the author's final selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Add, test, finalize dash and animation"
$ws.Range("A8").ClearContents() | Out-Null

$ws.Range("B18").Select() | Out-Null
